# Update capital structure database values for rows 2 and 3
# (Philippines / Reinsurance industry group) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 0.165                 # D2
$ws.Cells.Item(2, 7).Value = 0.07571644042232277   # G2
$ws.Cells.Item(2, 8).Value = 0.07571644042232277   # H2
$ws.Cells.Item(2, 9).Value = 0.02449078417551149   # I2
$ws.Cells.Item(2, 10).Value = 0.01224539208775574  # J2
$ws.Cells.Item(2, 11).Value = -0.836               # K2
$ws.Cells.Item(2, 12).Value = -0.01260935143288085 # L2
$ws.Cells.Item(2, 15).Value = -0                   # O2
$ws.Cells.Item(2, 18).Value = -0                   # R2
$ws.Cells.Item(2, 21).Value = 18                   # U2
$ws.Cells.Item(2, 22).Value = 0.6081081081081081   # V2
$ws.Cells.Item(2, 23).Value = -0.008148148148148147 # W2
$ws.Cells.Item(2, 24).Value = 0.09159116291327829  # X2
$ws.Cells.Item(2, 25).Value = -0.09973931106142644 # Y2
$ws.Cells.Item(2, 26).Value = 0.8575164056366941   # Z2
$ws.Cells.Item(2, 27).Value = 0.01050062460870432  # AA2
$ws.Cells.Item(2, 28).Value = 0.09147391432819912  # AB2
$ws.Cells.Item(2, 29).Value = -0.0809732897194948  # AC2
$ws.Cells.Item(2, 30).Value = 0.04                 # AD2
$ws.Cells.Item(2, 31).Value = 0.01630504581794165  # AE2
$ws.Cells.Item(2, 32).Value = 0.05630504581794165  # AF2
$ws.Cells.Item(2, 33).Value = -17.94369495418206   # AG2
$ws.Cells.Item(2, 34).Value = 0.001898586008302529 # AH2
$ws.Cells.Item(2, 35).Value = 0.0004872520438899182 # AI2
$ws.Cells.Item(2, 36).Value = -1.539398195538809   # AJ2
$ws.Cells.Item(2, 37).Value = -0.1839316786931884  # AK2
$ws.Cells.Item(2, 40).Value = 0.02371072910491998  # AN2
$ws.Cells.Item(2, 42).Value = -10.63645225499826   # AP2

# Row 3
$ws.Cells.Item(3, 4).Value = 0.165                 # D3
$ws.Cells.Item(3, 7).Value = 0.07571644042232277   # G3
$ws.Cells.Item(3, 8).Value = 0.07571644042232277   # H3
$ws.Cells.Item(3, 9).Value = 0.02449078417551149   # I3
$ws.Cells.Item(3, 10).Value = 0.01224539208775574  # J3
$ws.Cells.Item(3, 11).Value = -0.836               # K3
$ws.Cells.Item(3, 12).Value = -0.01260935143288085 # L3
$ws.Cells.Item(3, 15).Value = 0                    # O3
$ws.Cells.Item(3, 18).Value = 0                    # R3
$ws.Cells.Item(3, 21).Value = 18                   # U3
$ws.Cells.Item(3, 22).Value = 0.6081081081081081   # V3
$ws.Cells.Item(3, 23).Value = -0.008148148148148147 # W3
$ws.Cells.Item(3, 24).Value = 0.09159116291327829  # X3
$ws.Cells.Item(3, 25).Value = -0.09973931106142644 # Y3
$ws.Cells.Item(3, 26).Value = 0.8575164056366941   # Z3
$ws.Cells.Item(3, 27).Value = 0.01050062460870432  # AA3
$ws.Cells.Item(3, 28).Value = 0.09147391432819912  # AB3
$ws.Cells.Item(3, 29).Value = -0.0809732897194948  # AC3
$ws.Cells.Item(3, 30).Value = 0.04                 # AD3
$ws.Cells.Item(3, 31).Value = 0.01630504581794165  # AE3
$ws.Cells.Item(3, 32).Value = 0.05630504581794165  # AF3
$ws.Cells.Item(3, 33).Value = -17.94369495418206   # AG3
$ws.Cells.Item(3, 34).Value = 0.001898586008302529 # AH3
$ws.Cells.Item(3, 35).Value = 0.0004872520438899182 # AI3
$ws.Cells.Item(3, 36).Value = -1.539398195538809   # AJ3
$ws.Cells.Item(3, 37).Value = -0.1839316786931884  # AK3
$ws.Cells.Item(3, 40).Value = 0.02371072910491998  # AN3
$ws.Cells.Item(3, 42).Value = -10.63645225499826   # AP3
